$wb = $excel.ActiveWorkbook

# --- "About" sheet: bump the "last updated" date (2024-03-15 -> 2024-03-28) ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("C1").Value = (Get-Date -Year 2024 -Month 3 -Day 28 -Hour 0 -Minute 0 -Second 0)

# --- "RAF-capacity" sheet: hydrogen techs RAF raised from 0.3 to 1 ---
$wsCap = $wb.Worksheets.Item("RAF-capacity")
$wsCap.Range("B24").Value = 1
$wsCap.Range("B25").Value = 1

# RAF-capacity is now the active/selected tab
$wsCap.Activate()

# update the view for RAF-capacity: scrolled/zoomed while reviewing, selection on B25
$wsCap.Columns.Item(1).ColumnWidth = 28.166666666666668
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.Zoom = 80
$wsCap.Range("B25").Select()
